# Apply "Moved To Comp Bot" changes to the Mappings workbook.
# This adds five new CAN Talon mapping rows to the CANTalonSRX sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CANTalonSRX")

# New function names to append (rows 10-14), matching B column IDs 8-12.
# Set in the same order the original author typed them so that the shared
# string table is built up in the matching sequence.
$ws.Range("A10").Value = "ArmLift"
$ws.Range("A12").Value = "RobotLift1"
$ws.Range("A11").Value = "BallIntake"
$ws.Range("A13").Value = "RobotLift2"
$ws.Range("A14").Value = "RobotLiftExtend"

$ws.Range("B13").Value = 11
$ws.Range("B14").Value = 12

# Adjust column C width and selection as in the saved file.
# (ColumnWidth adds ~0.8333 padding when serialized to the XML "width"
# attribute, so back the value off to land exactly on 19.)
$ws.Columns.Item(3).ColumnWidth = 18.1666667
$ws.Range("C8").Select()

$wb.Save()
